$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# The table has a fixed row count (46) both before and after this edit: two
# rows get added and two get removed from the small block around rows 5-7,
# which nets out to zero, so instead of inserting/deleting rows we simply
# rewrite the text of every cell whose value actually changes. This also
# naturally handles the "collapse multi-run row down to one run" cases,
# since assigning Range.Text replaces all runs in that cell with a single run.

$t.Cell(1,1).Range.Text  = "0M"
$t.Cell(2,1).Range.Text  = "0M"
$t.Cell(3,1).Range.Text  = "0M"
$t.Cell(4,1).Range.Text  = "267"
$t.Cell(5,1).Range.Text  = "0.00002"
$t.Cell(6,1).Range.Text  = "0.00012"
$t.Cell(7,1).Range.Text  = "0.00004"
$t.Cell(12,1).Range.Text = "0.01034"
$t.Cell(44,1).Range.Text = "100"
$t.Cell(45,1).Range.Text = "0.01"
$t.Cell(46,1).Range.Text = "849"
